# Refresh the cryptocurrency price / 1h-volume snapshot (scraped data update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.639.91"
$ws.Range("E2").Value = "  -3.07%  "
$ws.Range("D3").Value = "2.084.70"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("B4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.64"
$ws.Range("B5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("B6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5152"
$ws.Range("B7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4396"
$ws.Range("B8").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("E8").Value = "  -2.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09184"
$ws.Range("B9").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.62"
$ws.Range("B10").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.174"
$ws.Range("B11").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.30"
$ws.Range("B12").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E12").Value = "  +3.75%  "
$ws.Range("D13").Value = "2.087.38"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.717"
$ws.Range("B14").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.119"
$ws.Range("B15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.32"
$ws.Range("B16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001158"
$ws.Range("B17").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.012"
$ws.Range("B18").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.88"
$ws.Range("B19").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").Value = "  +8.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06668"
$ws.Range("B20").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.009"
$ws.Range("B21").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.179"
$ws.Range("B22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("D23").Value = "29.749.52"
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.69"
$ws.Range("B24").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.312"
$ws.Range("B25").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = "  -3.04%  "
$ws.Range("D26").Value = "2.337.63"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.88"
$ws.Range("B27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.40"
$ws.Range("B28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.518"
$ws.Range("B29").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.88"
$ws.Range("B30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.142"
$ws.Range("B31").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "  -4.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1051"
$ws.Range("B32").Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("E32").Value = "  -1.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.618"
$ws.Range("B33").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.189"
$ws.Range("B34").Copy()
$ws.Range("D34").PasteSpecial(-4122)
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.952"
$ws.Range("B35").Copy()
$ws.Range("D35").PasteSpecial(-4122)
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.153"
$ws.Range("B36").Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("E36").Value = "  +4.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.12"
$ws.Range("B37").Copy()
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02573"
$ws.Range("B38").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06723"
$ws.Range("B39").Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2274"
$ws.Range("B40").Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.37"
$ws.Range("B41").Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("E41").Value = "  -1.92%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6848"
$ws.Range("B42").Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6631"
$ws.Range("B44").Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("E44").Value = "  +3.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.08"
$ws.Range("B45").Copy()
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("E45").Value = "  -5.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.296"
$ws.Range("B46").Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.626"
$ws.Range("B47").Copy()
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.215"
$ws.Range("B48").Copy()
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("E48").Value = "  -3.08%  "
$ws.Range("E49").Value = "  -6.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.73"
$ws.Range("B50").Copy()
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07103"
$ws.Range("B51").Copy()
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("E51").Value = "  -2.50%  "
$excel.CutCopyMode = 0
